# Updated cryptos list - price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.513.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.872.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  +0.72%  '
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4785'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3783'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07364'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9392'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("E11").Value = '  +5.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07850'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.881.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.446'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.590'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.69%  '
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008931'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.544.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.145'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.959'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("E26").Value = '  +2.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.021'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.012'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08934'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.334'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.217'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.97%  '
$ws.Range("E33").Value = '  +2.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7535'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.713'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02061'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.74%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.118'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05290'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.006'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5359'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.090'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.90%  '
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.457'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4827'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.27%  '
$ws.Range("E46").Value = '  +0.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.661'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06094'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9266'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.49%  '
